$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = '${CONCENTRADOR_UBI_DIRECCION}'
$ws.Range("F6").Value = '${EXTREMO_UBI_DIRECCION}'
[void]$ws.Range("F7").Select()
